$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5.623499999999999
$ws.Range("B3").Value = 5.823699999999992
$ws.Range("D4").Value = -6.803199999999995
$ws.Range("B5").Value = 4.742500000000001
$ws.Range("D6").Value = -8.703499999999993
$ws.Range("C7").Value = -12.6072
$ws.Range("A9").Value = -20.43219999999998
$ws.Range("C9").Value = -12.65460000000001
$ws.Range("D10").Value = -6.412399999999998
$ws.Range("B11").Value = 5.3958
$ws.Range("D11").Value = -8.324000000000005
$ws.Range("B12").Value = 5.500300000000001
$ws.Range("A13").Value = -21.88790000000002
$ws.Range("A16").Value = -19.94969999999999
$ws.Range("A18").Value = -22.72160000000001
$ws.Range("A20").Value = -22.04360000000002
$ws.Range("B21").Value = 5.790299999999994
$ws.Range("C21").Value = -11.8654
$ws.Range("D21").Value = -7.592700000000003
$ws.Range("D25").Value = -8.091099999999994
